# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet - this pushes the old "Over Due" column (N) to O and the old "#"
# column (P) to Q, inheriting column M's formatting for the new column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with R6 selected - this
# matches the workbook's saved view state (activeTab) and deselects the
# previously active "Transactions" sheet.
$ws.Activate()
$ws.Range("R6").Select()
